$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data
$ws.Range("A6").Value = 1551
$ws.Range("B6").Value = "DEPARTAMENTO DE TRÂNSITO DE MINAS GERAIS - DETRAN "

# Set column B width (91 characters, accounting for Excel's internal pixel rounding)
$ws.Columns.Item(2).ColumnWidth = 90.17

# Update the active selection cell shown in the saved view
$ws.Range("E7").Select()
